# Fill in column C on the "Marc" sheet with the BKS values found for the
# Golden instances (ρ = 10% results), matching the commit
# "working on reading Golden instances".
#
# The values are entered as literal text (not numbers) - we build them via a
# helper cell holding a text formula, copy it, and paste-special "values
# only" into the destination. That way the destination cell receives the
# already-typed string without Excel re-parsing it as a numeric literal
# (which would store it as a number and strip the "numbers as text"
# formatting the original workbook uses).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marc")

# row -> text value, in the same order the values were originally typed in
# (this also controls the order new entries land in the shared-strings
# table, so keep it as-is).
$entries = @(
    @{ Row = 9;  Text = "3500.7158" },
    @{ Row = 4;  Text = "1986.89" },
    @{ Row = 5;  Text = "2902.48" },
    @{ Row = 10; Text = "3215.310" },
    @{ Row = 7;  Text = "3075.87" },
    @{ Row = 8;  Text = "2251.90" },
    @{ Row = 3;  Text = "2094.59" },
    @{ Row = 6;  Text = "2730.05" }
)

$helper = $ws.Cells.Item(20, 5)

foreach ($entry in $entries) {
    $helper.Formula = '="' + $entry.Text + '"'
    $helper.Copy()
    $target = $ws.Cells.Item($entry.Row, 3)
    $target.PasteSpecial(-4163) # xlPasteValues
    $helper.ClearContents()
}

# Highlight the a-n20-c5.map row (row 6) BKS-at-10% value in red, the way
# the source workbook calls out that particular result.
$c6 = $ws.Range("C6")
$c6.Font.Color = 255
$c6.Font.Bold = $false
$c6.Font.Name = "Arial"

# The "Marc" sheet becomes the active tab/selection instead of "Golden".
$ws.Activate()
$ws.Range("C14").Select()
